# Add "Price" column (N) with per-row closing price data, formatted as currency.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("N1").Value = "Price"

# Price values for rows 7-36 (rows 2-6 intentionally left blank, matching source)
$prices = @{
    7  = 30.28
    8  = 32.05
    9  = 29.81
    10 = 41
    11 = 39.28
    12 = 37.85
    13 = 43.4
    14 = 45.2
    15 = 33.84
    16 = 38.62
    17 = 33.43
    18 = 26.9
    19 = 19.04
    20 = 20.49
    21 = 19.75
    22 = 7.79
    23 = 16.24
    24 = 10.5
    25 = 26.75
    26 = 23.54
    27 = 24.79
    28 = 32.91
    29 = 33.31
    30 = 37.45
    31 = 45.7
    32 = 39.18
    33 = 39.38
    34 = 37.45
    35 = 45.7
    36 = 39.18
}

foreach ($row in $prices.Keys) {
    $cell = $ws.Cells.Item($row, 14)   # column N = 14
    $cell.Value = $prices[$row]
    $cell.NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
}

# Touch N2:N6 so the column is fully formatted through the header block,
# leaving the values blank (matches the source workbook).
for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 14)
    $cell.Font.Name = "Calibri"
    $cell.Font.Color = 0
}
